$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7656.077
$ws.Range("I62").Value = 8471.843999999999
$ws.Range("J62").Value = 6350.85
$ws.Range("K62").Value = 8471.843999999999
$ws.Range("L62").Value = 6350.85
$ws.Range("M62").Value = -7847.843999999999
$ws.Range("N62").Value = -7598.85
$ws.Range("H65").Value = 7656.077
$ws.Range("I65").Value = 8471.843999999999
$ws.Range("J65").Value = 6350.85
$ws.Range("K65").Value = 42359.21999999999
$ws.Range("L65").Value = 31754.25
$ws.Range("M65").Value = -39239.21999999999
$ws.Range("N65").Value = -37994.25
$ws.Range("H74").Value = 6064.6577
$ws.Range("I74").Value = 3586.923
$ws.Range("K74").Value = 3586.923
$ws.Range("M74").Value = -2650.923
$ws.Range("H77").Value = 6064.6577
$ws.Range("I77").Value = 3586.923
$ws.Range("K77").Value = 17934.615
$ws.Range("M77").Value = -13254.615
$ws.Range("H112").Value = 1304.4324
$ws.Range("J112").Value = 1316.5714
$ws.Range("L112").Value = 3949.7142
$ws.Range("N112").Value = -6165.7142
$ws.Range("H138").Value = 1961.7537
$ws.Range("I138").Value = 1308.9143
$ws.Range("J138").Value = 2633.7942
$ws.Range("K138").Value = 3926.7429
$ws.Range("L138").Value = 7901.382599999999
$ws.Range("M138").Value = 1213.2571
$ws.Range("N138").Value = -18181.3826

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2909172.5
$ws.Range("I74").Value = 3473486
$ws.Range("K74").Value = 3473486
$ws.Range("M74").Value = -3472612
$ws.Range("H77").Value = 2909172.5
$ws.Range("I77").Value = 3473486
$ws.Range("K77").Value = 17367430
$ws.Range("M77").Value = -17363062
$ws.Range("H132").Value = 864794.25
$ws.Range("I132").Value = 932977.8
$ws.Range("K132").Value = 2798933.4
$ws.Range("M132").Value = -2796403.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 475.86667
$ws.Range("I22").Value = 475.86667
$ws.Range("K22").Value = 475.86667
$ws.Range("M22").Value = -302.86667
$ws.Range("H86").Value = 1237.8096
$ws.Range("I86").Value = 1067.8334
$ws.Range("J86").Value = 1464.4445
$ws.Range("K86").Value = 1067.8334
$ws.Range("L86").Value = 1464.4445
$ws.Range("M86").Value = 55.16660000000002
$ws.Range("N86").Value = -3710.4445
$ws.Range("H89").Value = 1237.8096
$ws.Range("I89").Value = 1067.8334
$ws.Range("J89").Value = 1464.4445
$ws.Range("K89").Value = 5339.166999999999
$ws.Range("L89").Value = 7322.2225
$ws.Range("M89").Value = 276.8330000000005
$ws.Range("N89").Value = -18554.2225
$ws.Range("H107").Value = 2694.3809
$ws.Range("I107").Value = 2412.9333
$ws.Range("J107").Value = 3398
$ws.Range("K107").Value = 2412.9333
$ws.Range("L107").Value = 3398
$ws.Range("M107").Value = -492.9333000000001
$ws.Range("N107").Value = -7238
$ws.Range("H134").Value = 1492987.5
$ws.Range("I134").Value = 2647923.2
$ws.Range("K134").Value = 7943769.600000001
$ws.Range("M134").Value = -7941234.600000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 75035.03999999999
$ws.Range("I31").Value = 106836.4
$ws.Range("J31").Value = 24822.37
$ws.Range("K31").Value = 106836.4
$ws.Range("L31").Value = 24822.37
$ws.Range("M31").Value = -106541.4
$ws.Range("N31").Value = -25412.37
$ws.Range("H34").Value = 75035.03999999999
$ws.Range("I34").Value = 106836.4
$ws.Range("J34").Value = 24822.37
$ws.Range("K34").Value = 106836.4
$ws.Range("L34").Value = 24822.37
$ws.Range("M34").Value = -106634.4
$ws.Range("N34").Value = -25226.37
$ws.Range("H58").Value = 318798.12
$ws.Range("I58").Value = 413150.84
$ws.Range("J58").Value = 4289.1113
$ws.Range("K58").Value = 413150.84
$ws.Range("L58").Value = 4289.1113
$ws.Range("M58").Value = -412947.84
$ws.Range("N58").Value = -4695.1113
$ws.Range("H86").Value = 3437.25
$ws.Range("I86").Value = 2900
$ws.Range("K86").Value = 2900
$ws.Range("M86").Value = -1777
$ws.Range("H89").Value = 3437.25
$ws.Range("I89").Value = 2900
$ws.Range("K89").Value = 14500
$ws.Range("M89").Value = -8884
$ws.Range("H99").Value = 3563.5862
$ws.Range("I99").Value = 3392.7778
$ws.Range("J99").Value = 3843.0908
$ws.Range("K99").Value = 3392.7778
$ws.Range("L99").Value = 3843.0908
$ws.Range("M99").Value = -1894.7778
$ws.Range("N99").Value = -6839.0908
$ws.Range("H126").Value = 3563.5862
$ws.Range("I126").Value = 3392.7778
$ws.Range("J126").Value = 3843.0908
$ws.Range("K126").Value = 10178.3334
$ws.Range("L126").Value = 11529.2724
$ws.Range("M126").Value = -7708.3334
$ws.Range("N126").Value = -16469.2724
$ws.Range("H132").Value = 12519784
$ws.Range("I132").Value = 33606.91
$ws.Range("J132").Value = 27780668
$ws.Range("K132").Value = 100820.73
$ws.Range("L132").Value = 83342004
$ws.Range("M132").Value = -98290.73000000001
$ws.Range("N132").Value = -83347064
$ws.Range("H136").Value = 318798.12
$ws.Range("I136").Value = 413150.84
$ws.Range("J136").Value = 4289.1113
$ws.Range("K136").Value = 1239452.52
$ws.Range("L136").Value = 12867.3339
$ws.Range("M136").Value = -1236902.52
$ws.Range("N136").Value = -17967.3339

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 73
$ws.Range("J12").Value = 69.59999999999999
$ws.Range("L12").Value = 208.8
$ws.Range("N12").Value = -554.8
$ws.Range("H140").Value = 2476.5938
$ws.Range("I140").Value = 1894.2693
$ws.Range("K140").Value = 5682.8079
$ws.Range("M140").Value = -502.8078999999998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 12678
$ws.Range("I9").Value = 250
$ws.Range("J9").Value = 20963.334
$ws.Range("K9").Value = 250
$ws.Range("L9").Value = 20963.334
$ws.Range("N9").Value = -21303.334
$ws.Range("M9").Value = -80
$ws.Range("H107").Value = 21638.588
$ws.Range("I107").Value = 25503.857
$ws.Range("J107").Value = 3600.6667
$ws.Range("K107").Value = 25503.857
$ws.Range("L107").Value = 3600.6667
$ws.Range("M107").Value = -23583.857
$ws.Range("N107").Value = -7440.6667
$ws.Range("H122").Value = 7247.75
$ws.Range("I122").Value = 2897.3
$ws.Range("J122").Value = 29000
$ws.Range("K122").Value = 8691.900000000001
$ws.Range("L122").Value = 87000
$ws.Range("M122").Value = -6241.900000000001
$ws.Range("N122").Value = -91900
$ws.Range("H132").Value = 1507400.1
$ws.Range("I132").Value = 1507400.1
$ws.Range("K132").Value = 4522200.300000001
$ws.Range("M132").Value = -4519670.300000001
$ws.Range("H140").Value = 108681.43
$ws.Range("J140").Value = 108681.43
$ws.Range("L140").Value = 108681.43
$ws.Range("N140").Value = -119041.43

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3551.3462
$ws.Range("I61").Value = 1602.8334
$ws.Range("J61").Value = 5221.5
$ws.Range("K61").Value = 1602.8334
$ws.Range("L61").Value = 5221.5
$ws.Range("M61").Value = -1400.8334
$ws.Range("N61").Value = -5625.5
$ws.Range("H69").Value = 70163
$ws.Range("J69").Value = 70163
$ws.Range("L69").Value = 70163
$ws.Range("N69").Value = -71785
$ws.Range("H72").Value = 70163
$ws.Range("J72").Value = 70163
$ws.Range("L72").Value = 210489
$ws.Range("N72").Value = -218601
$ws.Range("H82").Value = 989.75
$ws.Range("I82").Value = 735.2
$ws.Range("J82").Value = 1283.4615
$ws.Range("K82").Value = 735.2
$ws.Range("L82").Value = 1283.4615
$ws.Range("M82").Value = -374.2
$ws.Range("N82").Value = -2005.4615
$ws.Range("H85").Value = 989.75
$ws.Range("I85").Value = 735.2
$ws.Range("J85").Value = 1283.4615
$ws.Range("K85").Value = 735.2
$ws.Range("L85").Value = 1283.4615
$ws.Range("M85").Value = 512.8
$ws.Range("N85").Value = -3779.4615
$ws.Range("H113").Value = 3551.3462
$ws.Range("I113").Value = 1602.8334
$ws.Range("J113").Value = 5221.5
$ws.Range("K113").Value = 1602.8334
$ws.Range("L113").Value = 5221.5
$ws.Range("M113").Value = 567.1666
$ws.Range("N113").Value = -9561.5
$ws.Range("H136").Value = 27405.865
$ws.Range("I136").Value = 3109.0698
$ws.Range("K136").Value = 9327.2094
$ws.Range("M136").Value = -6777.2094

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4517
$ws.Range("I96").Value = 4201
$ws.Range("J96").Value = 4833
$ws.Range("K96").Value = 4201
$ws.Range("L96").Value = 4833
$ws.Range("M96").Value = -2828
$ws.Range("N96").Value = -7579
$ws.Range("H107").Value = 2536.9722
$ws.Range("I107").Value = 1485.963
$ws.Range("J107").Value = 5690
$ws.Range("K107").Value = 4457.889
$ws.Range("L107").Value = 17070
$ws.Range("M107").Value = -2537.889
$ws.Range("N107").Value = -20910
$ws.Range("H113").Value = 1267.0889
$ws.Range("I113").Value = 594.9355
$ws.Range("J113").Value = 2755.4285
$ws.Range("K113").Value = 1784.8065
$ws.Range("L113").Value = 8266.2855
$ws.Range("M113").Value = 385.1934999999999
$ws.Range("N113").Value = -12606.2855
$ws.Range("H122").Value = 1800.2916
$ws.Range("I122").Value = 1725.0233
$ws.Range("K122").Value = 5175.0699
$ws.Range("M122").Value = -2725.0699
$ws.Range("H126").Value = 3706.375
$ws.Range("I126").Value = 3678
$ws.Range("J126").Value = 3905
$ws.Range("K126").Value = 11034
$ws.Range("L126").Value = 11715
$ws.Range("M126").Value = -8564
$ws.Range("N126").Value = -16655
$ws.Range("H132").Value = 10068124
$ws.Range("I132").Value = 16772513
$ws.Range("K132").Value = 50317539
$ws.Range("M132").Value = -50315009
$ws.Range("H136").Value = 13583514
$ws.Range("I136").Value = 15832253
$ws.Range("K136").Value = 47496759
$ws.Range("M136").Value = -47494209
